$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# --- Prepare new rows 82:85 with the same cell formatting as row 81 ---
$ws.Range("A81:F81").Copy()
$ws.Range("A82:F82").PasteSpecial(-4122)
$ws.Range("A83:F85").PasteSpecial(-4122)

# Row 82 - Menu / Click menu
$ws.Range("A82").Value = "Menu"
$ws.Range("B82").Value = "Click menu"
$ws.Range("C82").Value = "A"
$ws.Range("E82").Value = "Daniel"
$ws.Range("F82").Value = 'Sound for "Click" on buttons for menu elements'
$ws.Range("G82").Value = "sound"

# Row 83 - Gameplay / Z-Virus-Cloud
$ws.Range("A83").Value = "Gameplay"
$ws.Range("B83").Value = "Z-Virus-Cloud"
$ws.Range("C83").Value = "A"
$ws.Range("E83").Value = "Daniel"
$ws.Range("F83").Value = "kleiner scalieren, damit er sich gehnau auf einem Tile befindet und nicht r$([char]0x00FC)ber geht"

# Row 84 - NPC / TS_Opa (comment column entered before the task column, matching author's edit order)
$ws.Range("F84").Value = "Farben eventuel anpassen um besser auf dem Level zu sichten"
$ws.Range("A84").Value = "NPC"
$ws.Range("B84").Value = "TS_Opa"
$ws.Range("C84").Value = "A"
$ws.Range("E84").Value = "Tobi"

# Row 85 - Gameplay / Out of cloud damage
$ws.Range("A85").Value = "Gameplay"
$ws.Range("B85").Value = "Out of cloud damage"
$ws.Range("C85").Value = "A"
$ws.Range("E85").Value = "Toma"
$ws.Range("F85").Value = "balancing of Healt loss out of the cloud"

# --- View / column layout changes ---
$ws.Columns.Item(6).ColumnWidth = 81.42578125
$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Range("F85").Select()

# --- Re-apply the AutoFilter over the grown range, filtered down to "Toma" ---
$ws.Range("A1:G85").AutoFilter(5, @("Toma"), 7)
